$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.514.26'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.741.60'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'323.10"
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = "'0.4507"
$ws.Range('D8').Value = "'0.3524"
$ws.Range('E8').Value = '  -2.09%  '
$ws.Range('D9').Value = "'0.07380"
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').Value = "'41.23"
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  -2.32%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('E13').Value = '  -1.22%  '
$ws.Range('D14').Value = "'5.906"
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').Value = "'7.056"
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '1.740.89'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = "'91.44"
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').Value = "'0.06350"
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = "'1.000"
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = "'16.77"
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('D22').Value = "'5.727"
$ws.Range('E22').Value = '  -2.50%  '
$ws.Range('D23').Value = '27.551.19'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').Value = "'2.097"
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').Value = "'162.41"
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').Value = '1.939.82'
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('D29').Value = "'124.86"
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('D30').Value = "'2.034"
$ws.Range('E30').Value = '  -4.46%  '
$ws.Range('D31').Value = "'1.046"
$ws.Range('E31').Value = '  -5.30%  '
$ws.Range('D32').Value = "'0.09065"
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = "'5.379"
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  -4.69%  '
$ws.Range('D37').Value = "'0.05981"
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = "'0.2058"
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = "'4.894"
$ws.Range('E39').Value = '  -0.91%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.6236"
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').Value = "'1.184"
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').Value = "'1.375"
$ws.Range('E42').Value = '  -0.67%  '
$ws.Range('D43').Value = "'7.693"
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('D44').Value = "'13.17"
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').Value = "'3.704"
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').Value = "'0.5792"
$ws.Range('D47').Value = "'121.93"
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = "'1.927"
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('E49').Value = '  +0.60%  '
$ws.Range('D50').Value = "'1.110"
$ws.Range('E50').Value = '  -4.83%  '
